$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.925.62'
$ws.Range('E2').Value = '  +1.73%  '

$ws.Range('D3').Value = '2.571.79'
$ws.Range('E3').Value = '  +3.15%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '302.50'
$ws.Range('E5').Value = '  +2.96%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '96.69'
$ws.Range('E6').Value = '  +5.05%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.575'
$ws.Range('E7').Value = '  +1.62%  '

$ws.Range('E8').Value = '  -0.03%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.550'
$ws.Range('E9').Value = '  +1.74%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '36.73'
$ws.Range('E10').Value = '  +2.91%  '

$ws.Range('E11').Value = '  +2.36%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '7.63'
$ws.Range('E12').Value = '  +0.77%  '

$ws.Range('E13').Value = '  +7.76%  '

$ws.Range('D14').Value = '2.582.84'
$ws.Range('E14').Value = '  +3.22%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.882'
$ws.Range('E15').Value = '  +3.42%  '

$ws.Range('E16').Value = '  +3.26%  '

$ws.Range('D17').Value = '42.963.14'
$ws.Range('E17').Value = '  +1.63%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '13.06'
$ws.Range('E18').Value = '  +7.09%  '

$ws.Range('E19').Value = '  +4.36%  '

$ws.Range('E20').Value = '  +3.49%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '72.07'
$ws.Range('E21').Value = '  +0.36%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '254.55'
$ws.Range('E22').Value = '  -0.43%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.97'
$ws.Range('E23').Value = '  +3.84%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.11'

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '28.52'
$ws.Range('E25').Value = '  +0.46%  '

$ws.Range('E26').Value = '  +0.00%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.23'
$ws.Range('E27').Value = '  +4.37%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '37.77'
$ws.Range('E28').Value = '  +3.90%  '

$ws.Range('E29').Value = '  -3.46%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.08'
$ws.Range('E30').Value = '  +2.68%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '155.69'
$ws.Range('E31').Value = '  +3.83%  '

$ws.Range('E32').Value = '  +1.47%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '2.76'
$ws.Range('E33').Value = '  +2.38%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.38'
$ws.Range('E34').Value = '  -1.04%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0811'
$ws.Range('E35').Value = '  +3.00%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '18.25'
$ws.Range('E36').Value = '  +12.08%  '

$ws.Range('E37').Value = '  +1.81%  '

$ws.Range('E38').Value = '  +2.03%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '23.65'
$ws.Range('E39').Value = '  +0.00%  '

$ws.Range('E40').Value = '  +1.16%  '

$ws.Range('E41').Value = '  +29.92%  '

$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0310'
$ws.Range('E42').Value = '  +1.76%  '

$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.87'
$ws.Range('E43').Value = '  +2.43%  '

$ws.Range('D44').Value = '2.062.39'
$ws.Range('E44').Value = '  +3.43%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.999'
$ws.Range('E45').Value = '  +0.11%  '

$ws.Range('E46').Value = '  +6.31%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '85.21'
$ws.Range('E47').Value = '  +0.66%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '77.07'
$ws.Range('E48').Value = '  +15.11%  '

$ws.Range('D49').Value = '2.821.85'
$ws.Range('E49').Value = '  +3.18%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '105.98'
$ws.Range('E50').Value = '  +4.50%  '

$ws.Range('E51').Value = '  +3.76%  '
